$d = $word.ActiveDocument

$d.Content.Find.Execute("61-15=", $true, $false, $false, $false, $false, $true, 1, $false, "0+57=", 2) | Out-Null
$d.Content.Find.Execute("11+47=", $true, $false, $false, $false, $false, $true, 1, $false, "85-59=", 2) | Out-Null
$d.Content.Find.Execute("21+59=", $true, $false, $false, $false, $false, $true, 1, $false, "5+68=", 2) | Out-Null
$d.Content.Find.Execute("85-49=", $true, $false, $false, $false, $false, $true, 1, $false, "58-24=", 2) | Out-Null
$d.Content.Find.Execute("46-5=", $true, $false, $false, $false, $false, $true, 1, $false, "52+1=", 2) | Out-Null
$d.Content.Find.Execute("49+50=", $true, $false, $false, $false, $false, $true, 1, $false, "96-79=", 2) | Out-Null
$d.Content.Find.Execute("95-15=", $true, $false, $false, $false, $false, $true, 1, $false, "70-41=", 2) | Out-Null
$d.Content.Find.Execute("13+44=", $true, $false, $false, $false, $false, $true, 1, $false, "12+29=", 2) | Out-Null
$d.Content.Find.Execute("89-56=", $true, $false, $false, $false, $false, $true, 1, $false, "39-11=", 2) | Out-Null
$d.Content.Find.Execute("6+22=", $true, $false, $false, $false, $false, $true, 1, $false, "90-68=", 2) | Out-Null
$d.Content.Find.Execute("51+17=", $true, $false, $false, $false, $false, $true, 1, $false, "3+28=", 2) | Out-Null
$d.Content.Find.Execute("15+37=", $true, $false, $false, $false, $false, $true, 1, $false, "20+49=", 2) | Out-Null
$d.Content.Find.Execute("64+19=", $true, $false, $false, $false, $false, $true, 1, $false, "57+38=", 2) | Out-Null
$d.Content.Find.Execute("96-18=", $true, $false, $false, $false, $false, $true, 1, $false, "5+51=", 2) | Out-Null
$d.Content.Find.Execute("76-33=", $true, $false, $false, $false, $false, $true, 1, $false, "82-30=", 2) | Out-Null
$d.Content.Find.Execute("74-54=", $true, $false, $false, $false, $false, $true, 1, $false, "32+63=", 2) | Out-Null
$d.Content.Find.Execute("60-33=", $true, $false, $false, $false, $false, $true, 1, $false, "30+67=", 2) | Out-Null
$d.Content.Find.Execute("82-5=", $true, $false, $false, $false, $false, $true, 1, $false, "96-11=", 2) | Out-Null
$d.Content.Find.Execute("38+3=", $true, $false, $false, $false, $false, $true, 1, $false, "62-60=", 2) | Out-Null
$d.Content.Find.Execute("22+66=", $true, $false, $false, $false, $false, $true, 1, $false, "79-11=", 2) | Out-Null
$d.Content.Find.Execute("62+26=", $true, $false, $false, $false, $false, $true, 1, $false, "38+60=", 2) | Out-Null
$d.Content.Find.Execute("76-65=", $true, $false, $false, $false, $false, $true, 1, $false, "94-15=", 2) | Out-Null
$d.Content.Find.Execute("26+58=", $true, $false, $false, $false, $false, $true, 1, $false, "90-50=", 2) | Out-Null
$d.Content.Find.Execute("18+13=", $true, $false, $false, $false, $false, $true, 1, $false, "23+51=", 2) | Out-Null
$d.Content.Find.Execute("34+38=", $true, $false, $false, $false, $false, $true, 1, $false, "78+9=", 2) | Out-Null
$d.Content.Find.Execute("73+0=", $true, $false, $false, $false, $false, $true, 1, $false, "14+82=", 2) | Out-Null
$d.Content.Find.Execute("1+85=", $true, $false, $false, $false, $false, $true, 1, $false, "56+13=", 2) | Out-Null
$d.Content.Find.Execute("79-20=", $true, $false, $false, $false, $false, $true, 1, $false, "40-18=", 2) | Out-Null
$d.Content.Find.Execute("42+2=", $true, $false, $false, $false, $false, $true, 1, $false, "83-58=", 2) | Out-Null
$d.Content.Find.Execute("33-6=", $true, $false, $false, $false, $false, $true, 1, $false, "25+39=", 2) | Out-Null
$d.Content.Find.Execute("54+37=", $true, $false, $false, $false, $false, $true, 1, $false, "6+28=", 2) | Out-Null
$d.Content.Find.Execute("94-76=", $true, $false, $false, $false, $false, $true, 1, $false, "36-20=", 2) | Out-Null
$d.Content.Find.Execute("69-39=", $true, $false, $false, $false, $false, $true, 1, $false, "7+10=", 2) | Out-Null
$d.Content.Find.Execute("62-45=", $true, $false, $false, $false, $false, $true, 1, $false, "19+63=", 2) | Out-Null
$d.Content.Find.Execute("85-63=", $true, $false, $false, $false, $false, $true, 1, $false, "15+69=", 2) | Out-Null
$d.Content.Find.Execute("80-36=", $true, $false, $false, $false, $false, $true, 1, $false, "62-7=", 2) | Out-Null
$d.Content.Find.Execute("40+22=", $true, $false, $false, $false, $false, $true, 1, $false, "3+53=", 2) | Out-Null
$d.Content.Find.Execute("70-65=", $true, $false, $false, $false, $false, $true, 1, $false, "92-45=", 2) | Out-Null
$d.Content.Find.Execute("6+48=", $true, $false, $false, $false, $false, $true, 1, $false, "36+16=", 2) | Out-Null
$d.Content.Find.Execute("3+39=", $true, $false, $false, $false, $false, $true, 1, $false, "29+37=", 2) | Out-Null
$d.Content.Find.Execute("48+20=", $true, $false, $false, $false, $false, $true, 1, $false, "66-46=", 2) | Out-Null
$d.Content.Find.Execute("3+2=", $true, $false, $false, $false, $false, $true, 1, $false, "81-23=", 2) | Out-Null
$d.Content.Find.Execute("70-42=", $true, $false, $false, $false, $false, $true, 1, $false, "6+81=", 2) | Out-Null
$d.Content.Find.Execute("84-82=", $true, $false, $false, $false, $false, $true, 1, $false, "76+23=", 2) | Out-Null
$d.Content.Find.Execute("51-21=", $true, $false, $false, $false, $false, $true, 1, $false, "43+13=", 2) | Out-Null
$d.Content.Find.Execute("76-75=", $true, $false, $false, $false, $false, $true, 1, $false, "33+52=", 2) | Out-Null
$d.Content.Find.Execute("5+1=", $true, $false, $false, $false, $false, $true, 1, $false, "83-72=", 2) | Out-Null
$d.Content.Find.Execute("45+14=", $true, $false, $false, $false, $false, $true, 1, $false, "63-10=", 2) | Out-Null
$d.Content.Find.Execute("13+65=", $true, $false, $false, $false, $false, $true, 1, $false, "21-16=", 2) | Out-Null
$d.Content.Find.Execute("13+80=", $true, $false, $false, $false, $false, $true, 1, $false, "67-43=", 2) | Out-Null
$d.Content.Find.Execute("66+5=", $true, $false, $false, $false, $false, $true, 1, $false, "94-12=", 2) | Out-Null
$d.Content.Find.Execute("17+23=", $true, $false, $false, $false, $false, $true, 1, $false, "52-4=", 2) | Out-Null
$d.Content.Find.Execute("7+5=", $true, $false, $false, $false, $false, $true, 1, $false, "29-3=", 2) | Out-Null
$d.Content.Find.Execute("36+41=", $true, $false, $false, $false, $false, $true, 1, $false, "37+5=", 2) | Out-Null
$d.Content.Find.Execute("44+10=", $true, $false, $false, $false, $false, $true, 1, $false, "99-9=", 2) | Out-Null
$d.Content.Find.Execute("49+22=", $true, $false, $false, $false, $false, $true, 1, $false, "78-11=", 2) | Out-Null
$d.Content.Find.Execute("40-31=", $true, $false, $false, $false, $false, $true, 1, $false, "70-0=", 2) | Out-Null
$d.Content.Find.Execute("11-6=", $true, $false, $false, $false, $false, $true, 1, $false, "11+30=", 2) | Out-Null
$d.Content.Find.Execute("89-40=", $true, $false, $false, $false, $false, $true, 1, $false, "1+49=", 2) | Out-Null
$d.Content.Find.Execute("64+6=", $true, $false, $false, $false, $false, $true, 1, $false, "6-3=", 2) | Out-Null
$d.Content.Find.Execute("19-17=", $true, $false, $false, $false, $false, $true, 1, $false, "31+12=", 2) | Out-Null
$d.Content.Find.Execute("78-37=", $true, $false, $false, $false, $false, $true, 1, $false, "66-53=", 2) | Out-Null
$d.Content.Find.Execute("16+19=", $true, $false, $false, $false, $false, $true, 1, $false, "36-18=", 2) | Out-Null
$d.Content.Find.Execute("47+33=", $true, $false, $false, $false, $false, $true, 1, $false, "87-67=", 2) | Out-Null
$d.Content.Find.Execute("58-53=", $true, $false, $false, $false, $false, $true, 1, $false, "37-0=", 2) | Out-Null
$d.Content.Find.Execute("96+1=", $true, $false, $false, $false, $false, $true, 1, $false, "97-55=", 2) | Out-Null
$d.Content.Find.Execute("30+29=", $true, $false, $false, $false, $false, $true, 1, $false, "11+3=", 2) | Out-Null
$d.Content.Find.Execute("59-8=", $true, $false, $false, $false, $false, $true, 1, $false, "73-66=", 2) | Out-Null
$d.Content.Find.Execute("68+11=", $true, $false, $false, $false, $false, $true, 1, $false, "98-59=", 2) | Out-Null
$d.Content.Find.Execute("38+14=", $true, $false, $false, $false, $false, $true, 1, $false, "43-19=", 2) | Out-Null
$d.Content.Find.Execute("37+6=", $true, $false, $false, $false, $false, $true, 1, $false, "29+13=", 2) | Out-Null
$d.Content.Find.Execute("91-43=", $true, $false, $false, $false, $false, $true, 1, $false, "59-46=", 2) | Out-Null
$d.Content.Find.Execute("54+21=", $true, $false, $false, $false, $false, $true, 1, $false, "52-47=", 2) | Out-Null
$d.Content.Find.Execute("54-16=", $true, $false, $false, $false, $false, $true, 1, $false, "10+89=", 2) | Out-Null
$d.Content.Find.Execute("26+7=", $true, $false, $false, $false, $false, $true, 1, $false, "84-15=", 2) | Out-Null
$d.Content.Find.Execute("55-30=", $true, $false, $false, $false, $false, $true, 1, $false, "25+8=", 2) | Out-Null
$d.Content.Find.Execute("35+47=", $true, $false, $false, $false, $false, $true, 1, $false, "20+2=", 2) | Out-Null
$d.Content.Find.Execute("82-47=", $true, $false, $false, $false, $false, $true, 1, $false, "13+25=", 2) | Out-Null
$d.Content.Find.Execute("70-8=", $true, $false, $false, $false, $false, $true, 1, $false, "41+28=", 2) | Out-Null
$d.Content.Find.Execute("44+52=", $true, $false, $false, $false, $false, $true, 1, $false, "88+7=", 2) | Out-Null
$d.Content.Find.Execute("90-64=", $true, $false, $false, $false, $false, $true, 1, $false, "17+77=", 2) | Out-Null
$d.Content.Find.Execute("71-58=", $true, $false, $false, $false, $false, $true, 1, $false, "87-14=", 2) | Out-Null
$d.Content.Find.Execute("57-28=", $true, $false, $false, $false, $false, $true, 1, $false, "92-29=", 2) | Out-Null
$d.Content.Find.Execute("15-4=", $true, $false, $false, $false, $false, $true, 1, $false, "13-0=", 2) | Out-Null
$d.Content.Find.Execute("38+5=", $true, $false, $false, $false, $false, $true, 1, $false, "91-18=", 2) | Out-Null
$d.Content.Find.Execute("62+23=", $true, $false, $false, $false, $false, $true, 1, $false, "57+17=", 2) | Out-Null
$d.Content.Find.Execute("54-39=", $true, $false, $false, $false, $false, $true, 1, $false, "39-20=", 2) | Out-Null
$d.Content.Find.Execute("38+23=", $true, $false, $false, $false, $false, $true, 1, $false, "62-27=", 2) | Out-Null
$d.Content.Find.Execute("84-69=", $true, $false, $false, $false, $false, $true, 1, $false, "16+50=", 2) | Out-Null
$d.Content.Find.Execute("76-28=", $true, $false, $false, $false, $false, $true, 1, $false, "11+70=", 2) | Out-Null
$d.Content.Find.Execute("10+83=", $true, $false, $false, $false, $false, $true, 1, $false, "22+22=", 2) | Out-Null
$d.Content.Find.Execute("33-11=", $true, $false, $false, $false, $false, $true, 1, $false, "49+49=", 2) | Out-Null
$d.Content.Find.Execute("43-5=", $true, $false, $false, $false, $false, $true, 1, $false, "81-80=", 2) | Out-Null
$d.Content.Find.Execute("54+34=", $true, $false, $false, $false, $false, $true, 1, $false, "61-0=", 2) | Out-Null
$d.Content.Find.Execute("34+59=", $true, $false, $false, $false, $false, $true, 1, $false, "24+26=", 2) | Out-Null
$d.Content.Find.Execute("31+68=", $true, $false, $false, $false, $false, $true, 1, $false, "86-23=", 2) | Out-Null
$d.Content.Find.Execute("44+53=", $true, $false, $false, $false, $false, $true, 1, $false, "91-30=", 2) | Out-Null
$d.Content.Find.Execute("65+30=", $true, $false, $false, $false, $false, $true, 1, $false, "72-66=", 2) | Out-Null
$d.Content.Find.Execute("13+62=", $true, $false, $false, $false, $false, $true, 1, $false, "22-5=", 2) | Out-Null
$d.Content.Find.Execute("64-2=", $true, $false, $false, $false, $false, $true, 1, $false, "1+9=", 2) | Out-Null
